$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "Execute" = "Y" markers for rows 2 and 4 (row 3 already had it)
$ws.Range("C2").Value = "Y"
$ws.Range("C4").Value = "Y"

# Match the final active selection in the saved file
$ws.Range("C4").Select()
